$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: notification / confirm-dialog strings.
$ws.Range("A16").Value = "lang_noti_header"
$ws.Range("B16").Value = "Thông Báo!"
$ws.Range("C16").Value = "Notification!"

$ws.Range("A17").Value = "lang_noti_login"
$ws.Range("C17").Value = "Login to continue"
$ws.Range("B17").Value = "Bạn cần đăng nhập để tiếp tục"

$ws.Range("A18").Value = "lang_confirm_ok"
$ws.Range("B18").Value = "Đồng ý"
$ws.Range("C18").Value = "Ok"

$ws.Range("A19").Value = "lang_confirm_cancel"
$ws.Range("B19").Value = "Hủy bỏ"
$ws.Range("C19").Value = "Cancel"

# Apply the same styles used by the other data rows (s="1" for column A,
# s="2" for columns B/C) to the newly added rows.
$ws.Range("A16:A19").Style = $ws.Range("A12").Style
$ws.Range("B16:C19").Style = $ws.Range("B12").Style

# Fix the typo'd key for the existing "student tracking" row (row 11).
# Writing the corrected string last gives it the final shared-string slot,
# dropping the old unused "lang_student_traking" entry automatically.
$ws.Range("A11").Value = "lang_student_tracking"

# Move the active selection to A11, matching the recorded view state.
$ws.Range("A11").Select()
